# Split the paragraph "On affiche le temps à la fin (lorsqu'on est mort!)"
# into two runs: "On " (unchanged) and "affiche le temps à la fin
# (lorsqu'on est mort!)" (highlighted in red), per the target diff.

$d = $word.ActiveDocument
$prefix = "On "
$found = $false

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    # Match on ASCII-safe substrings to sidestep any accented-character
    # round-tripping issues; this paragraph is unique in the document.
    if ($t.StartsWith("On affiche le temps") -and $t -like "*mort!)*") {
        $splitPoint = $p.Range.Start + $prefix.Length
        $endPoint = $p.Range.End - 1   # exclude the paragraph mark
        $r = $d.Range($splitPoint, $endPoint)
        $r.Font.HighlightColorIndex = 6   # wdRed
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Target paragraph ('On affiche le temps...') not found"
}
